# Apply rotation of values between rows 2, 3 and 5 for columns A, Q, R, AC
# Rotation (using original/before values):
#   row2 <- row5 (orig)
#   row3 <- row2 (orig)
#   row5 <- row3 (orig)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "Q", "R", "AC")

# Capture original values before any modification (use Value2 to avoid
# returning a boxed Variant wrapper from the interop layer)
$orig2 = @{}
$orig3 = @{}
$orig5 = @{}
foreach ($col in $cols) {
    $orig2[$col] = $ws.Range(($col + "2")).Value2
    $orig3[$col] = $ws.Range(($col + "3")).Value2
    $orig5[$col] = $ws.Range(($col + "5")).Value2
}

foreach ($col in $cols) {
    $ws.Range(($col + "2")).Value2 = $orig5[$col]
    $ws.Range(($col + "3")).Value2 = $orig2[$col]
    $ws.Range(($col + "5")).Value2 = $orig3[$col]
}
